$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 17, shifting rows 17-20 down to 18-21
$ws.Rows("17").Insert()

# Fill new row 17 with data
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C17").Value = "Arica y Parinacota"
$ws.Range("D17").Value = 44468
$ws.Range("E17").Value = 15
$ws.Range("F17").Value = 100112013
$ws.Range("G17").Value = "Alcachofa"
$ws.Range("H17").Value = "Argentina(o)"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 17000
$ws.Range("L17").Value = 18000
$ws.Range("M17").Value = 17500
$ws.Range("N17").Value = "$/caja 50 unidades"
$ws.Range("O17").Value = "Región de Coquimbo"
$ws.Range("P17").Value = 350
$ws.Range("Q17").Value = 50
$ws.Range("R17").Value = "Hortaliza"
